$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated Price (D) / Volume(1h) (E) values row by row, matching
# the latest cryptos-list scrape. Price cells whose new reading looks like a
# plain number are switched to Text format first so Excel keeps them as text
# strings (consistent with the rest of the sheet) instead of auto-converting
# them to numbers and silently dropping significant trailing zeros.

$ws.Range("D2").Value = "51.528.51"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "3.108.38"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "387.61"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.10"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.28"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0859"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "3.592.61"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.66"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.86"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "3.094.77"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.01"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.89"
$ws.Range("E18").Value = "  +3.30%  "
$ws.Range("D19").Value = "51.597.01"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  +7.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.56"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.21"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.62"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.11"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.59"
$ws.Range("E27").Value = "  +4.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.26"
$ws.Range("E28").Value = "  -3.50%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.45"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.77"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("E34").Value = "  +5.91%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.92"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.42"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.10"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.69"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  -2.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.80"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.21"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.51"
$ws.Range("E47").Value = "  +4.34%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "2.078.82"
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  +4.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.924"
$ws.Range("E51").Value = "  +18.04%  "
